$wb = $excel.ActiveWorkbook

# Add the new "Sprite Addresses" worksheet at the end of the tab order (after "Vera")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sprite Addresses"

# Header row
$ws.Cells.Item(1, 1).Value = "Index"
$ws.Cells.Item(1, 2).Value = "Value"

# 135 rows of (index, hex sprite address) pairs
  $ws.Cells.Item(2, 1).Value = 0
  $ws.Cells.Item(2, 2).Value = "0XEA00"
  $ws.Cells.Item(3, 1).Value = 1
  $ws.Cells.Item(3, 2).Value = "0XEC00"
  $ws.Cells.Item(4, 1).Value = 2
  $ws.Cells.Item(4, 2).Value = "0XEE00"
  $ws.Cells.Item(5, 1).Value = 3
  $ws.Cells.Item(5, 2).Value = "0XF000"
  $ws.Cells.Item(6, 1).Value = 4
  $ws.Cells.Item(6, 2).Value = "0XF200"
  $ws.Cells.Item(7, 1).Value = 5
  $ws.Cells.Item(7, 2).Value = "0XF400"
  $ws.Cells.Item(8, 1).Value = 6
  $ws.Cells.Item(8, 2).Value = "0XF600"
  $ws.Cells.Item(9, 1).Value = 7
  $ws.Cells.Item(9, 2).Value = "0XF800"
  $ws.Cells.Item(10, 1).Value = 8
  $ws.Cells.Item(10, 2).Value = "0XFA00"
  $ws.Cells.Item(11, 1).Value = 9
  $ws.Cells.Item(11, 2).Value = "0XFC00"
  $ws.Cells.Item(12, 1).Value = 10
  $ws.Cells.Item(12, 2).Value = "0XFE00"
  $ws.Cells.Item(13, 1).Value = 11
  $ws.Cells.Item(13, 2).Value = "0X10000"
  $ws.Cells.Item(14, 1).Value = 12
  $ws.Cells.Item(14, 2).Value = "0X10200"
  $ws.Cells.Item(15, 1).Value = 13
  $ws.Cells.Item(15, 2).Value = "0X10400"
  $ws.Cells.Item(16, 1).Value = 14
  $ws.Cells.Item(16, 2).Value = "0X10600"
  $ws.Cells.Item(17, 1).Value = 15
  $ws.Cells.Item(17, 2).Value = "0X10800"
  $ws.Cells.Item(18, 1).Value = 16
  $ws.Cells.Item(18, 2).Value = "0X10A00"
  $ws.Cells.Item(19, 1).Value = 17
  $ws.Cells.Item(19, 2).Value = "0X10C00"
  $ws.Cells.Item(20, 1).Value = 18
  $ws.Cells.Item(20, 2).Value = "0X10E00"
  $ws.Cells.Item(21, 1).Value = 19
  $ws.Cells.Item(21, 2).Value = "0X11000"
  $ws.Cells.Item(22, 1).Value = 20
  $ws.Cells.Item(22, 2).Value = "0X11200"
  $ws.Cells.Item(23, 1).Value = 21
  $ws.Cells.Item(23, 2).Value = "0X11400"
  $ws.Cells.Item(24, 1).Value = 22
  $ws.Cells.Item(24, 2).Value = "0X11600"
  $ws.Cells.Item(25, 1).Value = 23
  $ws.Cells.Item(25, 2).Value = "0X11800"
  $ws.Cells.Item(26, 1).Value = 24
  $ws.Cells.Item(26, 2).Value = "0X11A00"
  $ws.Cells.Item(27, 1).Value = 25
  $ws.Cells.Item(27, 2).Value = "0X11C00"
  $ws.Cells.Item(28, 1).Value = 26
  $ws.Cells.Item(28, 2).Value = "0X11E00"
  $ws.Cells.Item(29, 1).Value = 27
  $ws.Cells.Item(29, 2).Value = "0X12000"
  $ws.Cells.Item(30, 1).Value = 28
  $ws.Cells.Item(30, 2).Value = "0X12200"
  $ws.Cells.Item(31, 1).Value = 29
  $ws.Cells.Item(31, 2).Value = "0X12400"
  $ws.Cells.Item(32, 1).Value = 30
  $ws.Cells.Item(32, 2).Value = "0X12600"
  $ws.Cells.Item(33, 1).Value = 31
  $ws.Cells.Item(33, 2).Value = "0X12800"
  $ws.Cells.Item(34, 1).Value = 32
  $ws.Cells.Item(34, 2).Value = "0X12A00"
  $ws.Cells.Item(35, 1).Value = 33
  $ws.Cells.Item(35, 2).Value = "0X12C00"
  $ws.Cells.Item(36, 1).Value = 34
  $ws.Cells.Item(36, 2).Value = "0X12E00"
  $ws.Cells.Item(37, 1).Value = 35
  $ws.Cells.Item(37, 2).Value = "0X13000"
  $ws.Cells.Item(38, 1).Value = 36
  $ws.Cells.Item(38, 2).Value = "0X13200"
  $ws.Cells.Item(39, 1).Value = 37
  $ws.Cells.Item(39, 2).Value = "0X13400"
  $ws.Cells.Item(40, 1).Value = 38
  $ws.Cells.Item(40, 2).Value = "0X13600"
  $ws.Cells.Item(41, 1).Value = 39
  $ws.Cells.Item(41, 2).Value = "0X13800"
  $ws.Cells.Item(42, 1).Value = 40
  $ws.Cells.Item(42, 2).Value = "0X13A00"
  $ws.Cells.Item(43, 1).Value = 41
  $ws.Cells.Item(43, 2).Value = "0X13C00"
  $ws.Cells.Item(44, 1).Value = 42
  $ws.Cells.Item(44, 2).Value = "0X13E00"
  $ws.Cells.Item(45, 1).Value = 43
  $ws.Cells.Item(45, 2).Value = "0X14000"
  $ws.Cells.Item(46, 1).Value = 44
  $ws.Cells.Item(46, 2).Value = "0X14200"
  $ws.Cells.Item(47, 1).Value = 45
  $ws.Cells.Item(47, 2).Value = "0X14400"
  $ws.Cells.Item(48, 1).Value = 46
  $ws.Cells.Item(48, 2).Value = "0X14600"
  $ws.Cells.Item(49, 1).Value = 47
  $ws.Cells.Item(49, 2).Value = "0X14800"
  $ws.Cells.Item(50, 1).Value = 48
  $ws.Cells.Item(50, 2).Value = "0X14A00"
  $ws.Cells.Item(51, 1).Value = 49
  $ws.Cells.Item(51, 2).Value = "0X14C00"
  $ws.Cells.Item(52, 1).Value = 50
  $ws.Cells.Item(52, 2).Value = "0X14E00"
  $ws.Cells.Item(53, 1).Value = 51
  $ws.Cells.Item(53, 2).Value = "0X15000"
  $ws.Cells.Item(54, 1).Value = 52
  $ws.Cells.Item(54, 2).Value = "0X15200"
  $ws.Cells.Item(55, 1).Value = 53
  $ws.Cells.Item(55, 2).Value = "0X15400"
  $ws.Cells.Item(56, 1).Value = 54
  $ws.Cells.Item(56, 2).Value = "0X15600"
  $ws.Cells.Item(57, 1).Value = 55
  $ws.Cells.Item(57, 2).Value = "0X15800"
  $ws.Cells.Item(58, 1).Value = 56
  $ws.Cells.Item(58, 2).Value = "0X15A00"
  $ws.Cells.Item(59, 1).Value = 57
  $ws.Cells.Item(59, 2).Value = "0X15C00"
  $ws.Cells.Item(60, 1).Value = 58
  $ws.Cells.Item(60, 2).Value = "0X15E00"
  $ws.Cells.Item(61, 1).Value = 59
  $ws.Cells.Item(61, 2).Value = "0X16000"
  $ws.Cells.Item(62, 1).Value = 60
  $ws.Cells.Item(62, 2).Value = "0X16200"
  $ws.Cells.Item(63, 1).Value = 61
  $ws.Cells.Item(63, 2).Value = "0X16400"
  $ws.Cells.Item(64, 1).Value = 62
  $ws.Cells.Item(64, 2).Value = "0X16600"
  $ws.Cells.Item(65, 1).Value = 63
  $ws.Cells.Item(65, 2).Value = "0X16800"
  $ws.Cells.Item(66, 1).Value = 64
  $ws.Cells.Item(66, 2).Value = "0X16A00"
  $ws.Cells.Item(67, 1).Value = 65
  $ws.Cells.Item(67, 2).Value = "0X16C00"
  $ws.Cells.Item(68, 1).Value = 66
  $ws.Cells.Item(68, 2).Value = "0X16E00"
  $ws.Cells.Item(69, 1).Value = 67
  $ws.Cells.Item(69, 2).Value = "0X17000"
  $ws.Cells.Item(70, 1).Value = 68
  $ws.Cells.Item(70, 2).Value = "0X17200"
  $ws.Cells.Item(71, 1).Value = 69
  $ws.Cells.Item(71, 2).Value = "0X17400"
  $ws.Cells.Item(72, 1).Value = 70
  $ws.Cells.Item(72, 2).Value = "0X17600"
  $ws.Cells.Item(73, 1).Value = 71
  $ws.Cells.Item(73, 2).Value = "0X17800"
  $ws.Cells.Item(74, 1).Value = 72
  $ws.Cells.Item(74, 2).Value = "0X17A00"
  $ws.Cells.Item(75, 1).Value = 73
  $ws.Cells.Item(75, 2).Value = "0X17C00"
  $ws.Cells.Item(76, 1).Value = 74
  $ws.Cells.Item(76, 2).Value = "0X17E00"
  $ws.Cells.Item(77, 1).Value = 75
  $ws.Cells.Item(77, 2).Value = "0X18000"
  $ws.Cells.Item(78, 1).Value = 76
  $ws.Cells.Item(78, 2).Value = "0X18200"
  $ws.Cells.Item(79, 1).Value = 77
  $ws.Cells.Item(79, 2).Value = "0X18400"
  $ws.Cells.Item(80, 1).Value = 78
  $ws.Cells.Item(80, 2).Value = "0X18600"
  $ws.Cells.Item(81, 1).Value = 79
  $ws.Cells.Item(81, 2).Value = "0X18800"
  $ws.Cells.Item(82, 1).Value = 80
  $ws.Cells.Item(82, 2).Value = "0X18A00"
  $ws.Cells.Item(83, 1).Value = 81
  $ws.Cells.Item(83, 2).Value = "0X18C00"
  $ws.Cells.Item(84, 1).Value = 82
  $ws.Cells.Item(84, 2).Value = "0X18E00"
  $ws.Cells.Item(85, 1).Value = 83
  $ws.Cells.Item(85, 2).Value = "0X19000"
  $ws.Cells.Item(86, 1).Value = 84
  $ws.Cells.Item(86, 2).Value = "0X19200"
  $ws.Cells.Item(87, 1).Value = 85
  $ws.Cells.Item(87, 2).Value = "0X19400"
  $ws.Cells.Item(88, 1).Value = 86
  $ws.Cells.Item(88, 2).Value = "0X19600"
  $ws.Cells.Item(89, 1).Value = 87
  $ws.Cells.Item(89, 2).Value = "0X19800"
  $ws.Cells.Item(90, 1).Value = 88
  $ws.Cells.Item(90, 2).Value = "0X19A00"
  $ws.Cells.Item(91, 1).Value = 89
  $ws.Cells.Item(91, 2).Value = "0X19C00"
  $ws.Cells.Item(92, 1).Value = 90
  $ws.Cells.Item(92, 2).Value = "0X19E00"
  $ws.Cells.Item(93, 1).Value = 91
  $ws.Cells.Item(93, 2).Value = "0X1A000"
  $ws.Cells.Item(94, 1).Value = 92
  $ws.Cells.Item(94, 2).Value = "0X1A200"
  $ws.Cells.Item(95, 1).Value = 93
  $ws.Cells.Item(95, 2).Value = "0X1A400"
  $ws.Cells.Item(96, 1).Value = 94
  $ws.Cells.Item(96, 2).Value = "0X1A600"
  $ws.Cells.Item(97, 1).Value = 95
  $ws.Cells.Item(97, 2).Value = "0X1A800"
  $ws.Cells.Item(98, 1).Value = 96
  $ws.Cells.Item(98, 2).Value = "0X1AA00"
  $ws.Cells.Item(99, 1).Value = 97
  $ws.Cells.Item(99, 2).Value = "0X1AC00"
  $ws.Cells.Item(100, 1).Value = 98
  $ws.Cells.Item(100, 2).Value = "0X1AE00"
  $ws.Cells.Item(101, 1).Value = 99
  $ws.Cells.Item(101, 2).Value = "0X1B000"
  $ws.Cells.Item(102, 1).Value = 100
  $ws.Cells.Item(102, 2).Value = "0X1B200"
  $ws.Cells.Item(103, 1).Value = 101
  $ws.Cells.Item(103, 2).Value = "0X1B400"
  $ws.Cells.Item(104, 1).Value = 102
  $ws.Cells.Item(104, 2).Value = "0X1B600"
  $ws.Cells.Item(105, 1).Value = 103
  $ws.Cells.Item(105, 2).Value = "0X1B800"
  $ws.Cells.Item(106, 1).Value = 104
  $ws.Cells.Item(106, 2).Value = "0X1BA00"
  $ws.Cells.Item(107, 1).Value = 105
  $ws.Cells.Item(107, 2).Value = "0X1BC00"
  $ws.Cells.Item(108, 1).Value = 106
  $ws.Cells.Item(108, 2).Value = "0X1BE00"
  $ws.Cells.Item(109, 1).Value = 107
  $ws.Cells.Item(109, 2).Value = "0X1C000"
  $ws.Cells.Item(110, 1).Value = 108
  $ws.Cells.Item(110, 2).Value = "0X1C200"
  $ws.Cells.Item(111, 1).Value = 109
  $ws.Cells.Item(111, 2).Value = "0X1C400"
  $ws.Cells.Item(112, 1).Value = 110
  $ws.Cells.Item(112, 2).Value = "0X1C600"
  $ws.Cells.Item(113, 1).Value = 111
  $ws.Cells.Item(113, 2).Value = "0X1C800"
  $ws.Cells.Item(114, 1).Value = 112
  $ws.Cells.Item(114, 2).Value = "0X1CA00"
  $ws.Cells.Item(115, 1).Value = 113
  $ws.Cells.Item(115, 2).Value = "0X1CC00"
  $ws.Cells.Item(116, 1).Value = 114
  $ws.Cells.Item(116, 2).Value = "0X1CE00"
  $ws.Cells.Item(117, 1).Value = 115
  $ws.Cells.Item(117, 2).Value = "0X1D000"
  $ws.Cells.Item(118, 1).Value = 116
  $ws.Cells.Item(118, 2).Value = "0X1D200"
  $ws.Cells.Item(119, 1).Value = 117
  $ws.Cells.Item(119, 2).Value = "0X1D400"
  $ws.Cells.Item(120, 1).Value = 118
  $ws.Cells.Item(120, 2).Value = "0X1D600"
  $ws.Cells.Item(121, 1).Value = 119
  $ws.Cells.Item(121, 2).Value = "0X1D800"
  $ws.Cells.Item(122, 1).Value = 120
  $ws.Cells.Item(122, 2).Value = "0X1DA00"
  $ws.Cells.Item(123, 1).Value = 121
  $ws.Cells.Item(123, 2).Value = "0X1DC00"
  $ws.Cells.Item(124, 1).Value = 122
  $ws.Cells.Item(124, 2).Value = "0X1DE00"
  $ws.Cells.Item(125, 1).Value = 123
  $ws.Cells.Item(125, 2).Value = "0X1E000"
  $ws.Cells.Item(126, 1).Value = 124
  $ws.Cells.Item(126, 2).Value = "0X1E200"
  $ws.Cells.Item(127, 1).Value = 125
  $ws.Cells.Item(127, 2).Value = "0X1E400"
  $ws.Cells.Item(128, 1).Value = 126
  $ws.Cells.Item(128, 2).Value = "0X1E600"
  $ws.Cells.Item(129, 1).Value = 127
  $ws.Cells.Item(129, 2).Value = "0X1E800"
  $ws.Cells.Item(130, 1).Value = 128
  $ws.Cells.Item(130, 2).Value = "0X1EA00"
  $ws.Cells.Item(131, 1).Value = 129
  $ws.Cells.Item(131, 2).Value = "0X1EC00"
  $ws.Cells.Item(132, 1).Value = 130
  $ws.Cells.Item(132, 2).Value = "0X1EE00"
  $ws.Cells.Item(133, 1).Value = 131
  $ws.Cells.Item(133, 2).Value = "0X1F000"
  $ws.Cells.Item(134, 1).Value = 132
  $ws.Cells.Item(134, 2).Value = "0X1F200"
  $ws.Cells.Item(135, 1).Value = 133
  $ws.Cells.Item(135, 2).Value = "0X1F400"
  $ws.Cells.Item(136, 1).Value = 134
  $ws.Cells.Item(136, 2).Value = "0X1F600"

# Wide descriptive note in column D, added last so it lands at the end of the shared-string table
$ws.Cells.Item(1, 4).Value = "All possible sprite addresses and their index in _bESpriteAllocTable"

# Column D is widened to fit the note
$ws.Columns.Item(4).ColumnWidth = 119.1796875

# Select D1 so it is the active cell on this (now active/selected) sheet
$ws.Range("D1").Select()
